# Append a new observation record as row 3 on the Artfynd sheet.
# (Matches the added <row r="3"> in the target OOXML diff; dimension
# grows from A1:AY2 to A1:AY3 automatically once these cells are set.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 111908957
$ws.Range("B3").Value = 90687
$ws.Range("C3").Value = "Ovaliderad"
$ws.Range("D3").Value = "LC"
$ws.Range("E3").Value = 5964
$ws.Range("F3").Value = "Fjällig taggsvamp s.str."
$ws.Range("G3").Value = "Sarcodon imbricatus s.str."
$ws.Range("H3").Value = "(L.:Fr.) P.Karst."
# Leading "'" forces text so "7" is stored as the string "7", not number 7.
$ws.Range("I3").Value = "'7"
$ws.Range("J3").Value = "fruktkroppar"
# Leading "'" with nothing after it yields a present-but-empty text cell
# (plain "" would just clear/omit the cell instead of creating it).
$ws.Range("K3").Value = "'"
$ws.Range("N3").Value = "'"
$ws.Range("P3").Value = "Prästgården, Upl"
$ws.Range("Q3").Value = 655897
$ws.Range("R3").Value = 6675360
$ws.Range("S3").Value = 4
$ws.Range("T3").Value = "Uppsala"
$ws.Range("U3").Value = "Östhammar"
$ws.Range("V3").Value = "Uppland"
$ws.Range("W3").Value = "Dannemora"
# Leading "'" keeps these date-looking strings as literal text instead of
# being parsed into date serial numbers.
$ws.Range("Y3").Value = "'2023-09-05"
$ws.Range("Z3").Value = "13:32"
$ws.Range("AA3").Value = "'2023-09-05"
$ws.Range("AB3").Value = "13:32"
$ws.Range("AD3").Value = $False
$ws.Range("AE3").Value = $False
$ws.Range("AG3").Value = $False
$ws.Range("AT3").Value = "'"
$ws.Range("AW3").Value = "Annika Rastén"
$ws.Range("AX3").Value = "Annika Rastén"
$ws.Range("AY3").Value = "'"
